# LandKreisKode.xlsx edit:
#  - Six "Status" (column G) corrections on sheet "Tabelle1"
#  - Scroll the view down (topLeftCell A13 -> A37) and move the
#    selection / active cell from G38 to G71.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column G value corrections -----------------------------------------
$ws.Range("G34").Value = -1   # was -2
$ws.Range("G37").Value = -1   # was -2
$ws.Range("G39").Value = 0    # was -1
$ws.Range("G54").Value = -1   # was 0
$ws.Range("G70").Value = -1   # was 0
$ws.Range("G72").Value = -1   # was 0

# --- View state: scroll + selection -------------------------------------
# Move the window's visible top-left corner to row 37 (column A) ...
$win = $excel.ActiveWindow
$win.ScrollRow = 37
$win.ScrollColumn = 1

# ... and set the active cell / selection to G71 (was G38).
$ws.Range("G71").Select() | Out-Null
